$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("U-R-B")
Write-Output $ws.Name
